# "Reading data from Excel class 32"
# Add a new "Employee" sheet, and turn Sheet2 into a small username/password
# data table (with mailto: hyperlinks on the email/username column) used by
# the Selenium/SDET "reading data from Excel" exercises.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1. Populate Sheet2 with the Username / Password table.
# ---------------------------------------------------------------------
$ws2.Range("A1").Value = "Username"
$ws2.Range("B1").Value = "Password"

$ws2.Range("A2").Value = "asghar@gmail.com"
$ws2.Range("B2").Value = "123SKBDSSD+_#_$"

$ws2.Range("A3").Value = "asghar@gmail"
$ws2.Range("B3").Value = "123SKBDSSD+_#_$"

$ws2.Range("A4").Value = "asghar@gmail."
$ws2.Range("B4").Value = "123SKBDSSD+_#_$"

$ws2.Range("A5").Value = "asghargmail.com"
$ws2.Range("B5").Value = "123SKBDSSD+_#_$"

$ws2.Range("A6").Value = "asghargmail.com"
$ws2.Range("B6").Value = "123SKBDSSD"

$ws2.Range("A7").Value = "asghargmail.com"
$ws2.Range("B7").Value = "123SKBDSSD"

# ---------------------------------------------------------------------
# 2. Hyperlink the Username column (mailto: links).  For rows 2-4 the
#    link text matches the cell text exactly, so no "display" override is
#    written.  For rows 5-7 the cell text was edited after the hyperlink
#    was auto-created, so the original typed address is kept as the
#    hyperlink's TextToDisplay (-> OOXML "display" attribute) while the
#    cell itself shows the edited text.
# ---------------------------------------------------------------------
$ws2.Hyperlinks.Add($ws2.Range("A2"), "mailto:asghar@gmail.com") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "mailto:asghar@gmail") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "mailto:asghar@gmail.") | Out-Null

$ws2.Hyperlinks.Add($ws2.Range("A5"), "mailto:asghargmail.com", [Type]::Missing, [Type]::Missing, "asghar@gmail.") | Out-Null
$ws2.Range("A5").Value = "asghargmail.com"

$ws2.Hyperlinks.Add($ws2.Range("A6"), "mailto:asghargmail.com", [Type]::Missing, [Type]::Missing, "asghar@gmail.") | Out-Null
$ws2.Range("A6").Value = "asghargmail.com"

$ws2.Hyperlinks.Add($ws2.Range("A7"), "mailto:asghargmail.com", [Type]::Missing, [Type]::Missing, "asghar@gmail.") | Out-Null
$ws2.Range("A7").Value = "asghargmail.com"

# ---------------------------------------------------------------------
# 3. Column widths on Sheet2.
# ---------------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 29.14
$ws2.Columns.Item(2).ColumnWidth = 20.75

# ---------------------------------------------------------------------
# 4. Add the new (blank) "Employee" sheet after Sheet2.
# ---------------------------------------------------------------------
$wsEmployee = $wb.Worksheets.Add([Type]::Missing, $ws2)
$wsEmployee.Name = "Employee"

# ---------------------------------------------------------------------
# 5. Selection / view state: Sheet1's cursor moved to B24 (no longer the
#    active tab), Sheet2 becomes the active tab at 132% zoom with the
#    cursor on A10.
# ---------------------------------------------------------------------
$ws1.Range("B24").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("A10").Select() | Out-Null
$excel.ActiveWindow.Zoom = 132
